$d = $word.ActiveDocument
$t = $d.Tables(1)
$cell = $t.Cell(11, 2)

# Clear existing content; leaves a single empty paragraph
$cell.Range.Delete()
$cell.Range.Paragraphs(1).Alignment = 0

# --- paragraph 1: 11/12 – MRSA SCREEN – **Negative** ---
$cell.Range.Paragraphs(1).Range.InsertParagraphAfter()
$p1 = $cell.Range.Paragraphs(2).Range
$p1.InsertAfter("11/12 – MRSA SCREEN – **Negative**")
$p1 = $cell.Range.Paragraphs(2).Range
$p1.Font.Name = "Times New Roman"
$p1.Font.Color = 16711680
$p1.Font.Size = 10

# --- paragraph 2:   Summary: No MRSA isolated. ---
$cell.Range.Paragraphs(2).Range.InsertParagraphAfter()
$p2 = $cell.Range.Paragraphs(3).Range
$p2.InsertAfter("  Summary: No MRSA isolated.")
$p2 = $cell.Range.Paragraphs(3).Range
$p2.Font.Name = "Times New Roman"
$p2.Font.Color = 16711680
$p2.Font.Size = 10

# --- paragraph 3: 11/12 – CPE SCREEN – **Negative** ---
$cell.Range.Paragraphs(3).Range.InsertParagraphAfter()
$p3 = $cell.Range.Paragraphs(4).Range
$p3.InsertAfter("11/12 – CPE SCREEN – **Negative**")
$p3 = $cell.Range.Paragraphs(4).Range
$p3.Font.Name = "Times New Roman"
$p3.Font.Color = 16711680
$p3.Font.Size = 10

# --- paragraph 4:   Summary: No CPE organisms detected. ---
$cell.Range.Paragraphs(4).Range.InsertParagraphAfter()
$p4 = $cell.Range.Paragraphs(5).Range
$p4.InsertAfter("  Summary: No CPE organisms detected.")
$p4 = $cell.Range.Paragraphs(5).Range
$p4.Font.Name = "Times New Roman"
$p4.Font.Color = 16711680
$p4.Font.Size = 10

# --- paragraph 5: 10/12 – BLC – PICC LINE NO GROWTH AFTER 5 DAYS ---
$cell.Range.Paragraphs(5).Range.InsertParagraphAfter()
$p5 = $cell.Range.Paragraphs(6).Range
$p5.InsertAfter("10/12 – BLC – PICC LINE NO GROWTH AFTER 5 DAYS")
$p5 = $cell.Range.Paragraphs(6).Range
$p5.Font.Name = "Times New Roman"
$p5.Font.Color = 16711680
$p5.Font.Size = 10

# --- paragraph 6: 09/12 – LINE TIPS CULTURE – **No clear Result** ---
$cell.Range.Paragraphs(6).Range.InsertParagraphAfter()
$p6 = $cell.Range.Paragraphs(7).Range
$p6.InsertAfter("09/12 – LINE TIPS CULTURE – **No clear Result**")
$p6 = $cell.Range.Paragraphs(7).Range
$p6.Font.Name = "Times New Roman"
$p6.Font.Color = 16711680
$p6.Font.Size = 10

# --- paragraph 7:   Summary: Criteria for culturing IV line tips bas ---
$cell.Range.Paragraphs(7).Range.InsertParagraphAfter()
$p7 = $cell.Range.Paragraphs(8).Range
$p7.InsertAfter("  Summary: Criteria for culturing IV line tips based on blood cultures.")
$p7 = $cell.Range.Paragraphs(8).Range
$p7.Font.Name = "Times New Roman"
$p7.Font.Color = 16711680
$p7.Font.Size = 10

# --- paragraph 8: 08/12 – TB CULTURE – EARLY MORNING URINE (EMU) MYC ---
$cell.Range.Paragraphs(8).Range.InsertParagraphAfter()
$p8 = $cell.Range.Paragraphs(9).Range
$p8.InsertAfter("08/12 – TB CULTURE – EARLY MORNING URINE (EMU) MYCOBACTERIAL CULTURE: NO GROWTH AFTER 8 WEEKS")
$p8 = $cell.Range.Paragraphs(9).Range
$p8.Font.Name = "Times New Roman"
$p8.Font.Color = 16711680
$p8.Font.Size = 10

# --- paragraph 9: 08/12 – TB CULTURE – LYMPH NODE MYCOBACTERIAL CULT ---
$cell.Range.Paragraphs(9).Range.InsertParagraphAfter()
$p9 = $cell.Range.Paragraphs(10).Range
$p9.InsertAfter("08/12 – TB CULTURE – LYMPH NODE MYCOBACTERIAL CULTURE: NO GROWTH AFTER 8 WEEKS")
$p9 = $cell.Range.Paragraphs(10).Range
$p9.Font.Name = "Times New Roman"
$p9.Font.Color = 16711680
$p9.Font.Size = 10

# --- paragraph 10: 08/12 – NOCARDIA CULTURE – **Negative** ---
$cell.Range.Paragraphs(10).Range.InsertParagraphAfter()
$p10 = $cell.Range.Paragraphs(11).Range
$p10.InsertAfter("08/12 – NOCARDIA CULTURE – **Negative**")
$p10 = $cell.Range.Paragraphs(11).Range
$p10.Font.Name = "Times New Roman"
$p10.Font.Color = 16711680
$p10.Font.Size = 10

# --- paragraph 11:   Summary: Culture negative for Nocardia. ---
$cell.Range.Paragraphs(11).Range.InsertParagraphAfter()
$p11 = $cell.Range.Paragraphs(12).Range
$p11.InsertAfter("  Summary: Culture negative for Nocardia.")
$p11 = $cell.Range.Paragraphs(12).Range
$p11.Font.Name = "Times New Roman"
$p11.Font.Color = 16711680
$p11.Font.Size = 10

# --- paragraph 12: 08/12 – REF POSACONAZOLE LEVEL – **No clear Result ---
$cell.Range.Paragraphs(12).Range.InsertParagraphAfter()
$p12 = $cell.Range.Paragraphs(13).Range
$p12.InsertAfter("08/12 – REF POSACONAZOLE LEVEL – **No clear Result** +")
$p12 = $cell.Range.Paragraphs(13).Range
$p12.Font.Name = "Times New Roman"
$p12.Font.Color = 16711680
$p12.Font.Size = 10

# --- paragraph 13:   Summary: No result provided for Posaconazole Lev ---
$cell.Range.Paragraphs(13).Range.InsertParagraphAfter()
$p13 = $cell.Range.Paragraphs(14).Range
$p13.InsertAfter("  Summary: No result provided for Posaconazole Level.")
$p13 = $cell.Range.Paragraphs(14).Range
$p13.Font.Name = "Times New Roman"
$p13.Font.Color = 16711680
$p13.Font.Size = 10

# --- paragraph 14: 08/12 – ANAEROBIC CULTURE – **Negative** ---
$cell.Range.Paragraphs(14).Range.InsertParagraphAfter()
$p14 = $cell.Range.Paragraphs(15).Range
$p14.InsertAfter("08/12 – ANAEROBIC CULTURE – **Negative**")
$p14 = $cell.Range.Paragraphs(15).Range
$p14.Font.Name = "Times New Roman"
$p14.Font.Color = 16711680
$p14.Font.Size = 10

# --- paragraph 15:   Summary: No anaerobes isolated. ---
$cell.Range.Paragraphs(15).Range.InsertParagraphAfter()
$p15 = $cell.Range.Paragraphs(16).Range
$p15.InsertAfter("  Summary: No anaerobes isolated.")
$p15 = $cell.Range.Paragraphs(16).Range
$p15.Font.Name = "Times New Roman"
$p15.Font.Color = 16711680
$p15.Font.Size = 10

# --- paragraph 16: 08/12 – PCR FOR MTB COMPLEX – **Positive** ---
$cell.Range.Paragraphs(16).Range.InsertParagraphAfter()
$p16 = $cell.Range.Paragraphs(17).Range
$p16.InsertAfter("08/12 – PCR FOR MTB COMPLEX – **Positive**")
$p16 = $cell.Range.Paragraphs(17).Range
$p16.Font.Name = "Times New Roman"
$p16.Font.Color = 16711680
$p16.Font.Size = 10

# --- paragraph 17:   Summary: MTB complex detected, no Rifampicin res ---
$cell.Range.Paragraphs(17).Range.InsertParagraphAfter()
$p17 = $cell.Range.Paragraphs(18).Range
$p17.InsertAfter("  Summary: MTB complex detected, no Rifampicin resistance.")
$p17 = $cell.Range.Paragraphs(18).Range
$p17.Font.Name = "Times New Roman"
$p17.Font.Color = 16711680
$p17.Font.Size = 10

# --- paragraph 18: 08/12 – MICROSCOPY – **Negative** ---
$cell.Range.Paragraphs(18).Range.InsertParagraphAfter()
$p18 = $cell.Range.Paragraphs(19).Range
$p18.InsertAfter("08/12 – MICROSCOPY – **Negative**")
$p18 = $cell.Range.Paragraphs(19).Range
$p18.Font.Name = "Times New Roman"
$p18.Font.Color = 16711680
$p18.Font.Size = 10

# --- paragraph 19:   Summary: No cells or organisms seen. ---
$cell.Range.Paragraphs(19).Range.InsertParagraphAfter()
$p19 = $cell.Range.Paragraphs(20).Range
$p19.InsertAfter("  Summary: No cells or organisms seen.")
$p19 = $cell.Range.Paragraphs(20).Range
$p19.Font.Name = "Times New Roman"
$p19.Font.Color = 16711680
$p19.Font.Size = 10

# --- paragraph 20: 08/12 – FUNGUS CULTURE – **Negative** ---
$cell.Range.Paragraphs(20).Range.InsertParagraphAfter()
$p20 = $cell.Range.Paragraphs(21).Range
$p20.InsertAfter("08/12 – FUNGUS CULTURE – **Negative**")
$p20 = $cell.Range.Paragraphs(21).Range
$p20.Font.Name = "Times New Roman"
$p20.Font.Color = 16711680
$p20.Font.Size = 10

# --- paragraph 21:   Summary: Fungal cultures negative. ---
$cell.Range.Paragraphs(21).Range.InsertParagraphAfter()
$p21 = $cell.Range.Paragraphs(22).Range
$p21.InsertAfter("  Summary: Fungal cultures negative.")
$p21 = $cell.Range.Paragraphs(22).Range
$p21.Font.Name = "Times New Roman"
$p21.Font.Color = 16711680
$p21.Font.Size = 10

# --- paragraph 22: 08/12 – PUS MICRO / CULTURE – LYMPH NODE NO GROWTH ---
$cell.Range.Paragraphs(22).Range.InsertParagraphAfter()
$p22 = $cell.Range.Paragraphs(23).Range
$p22.InsertAfter("08/12 – PUS MICRO / CULTURE – LYMPH NODE NO GROWTH")
$p22 = $cell.Range.Paragraphs(23).Range
$p22.Font.Name = "Times New Roman"
$p22.Font.Color = 16711680
$p22.Font.Size = 10

# --- paragraph 23: 08/12 – TB MICROSCOPY – **Negative** ---
$cell.Range.Paragraphs(23).Range.InsertParagraphAfter()
$p23 = $cell.Range.Paragraphs(24).Range
$p23.InsertAfter("08/12 – TB MICROSCOPY – **Negative**")
$p23 = $cell.Range.Paragraphs(24).Range
$p23.Font.Name = "Times New Roman"
$p23.Font.Color = 16711680
$p23.Font.Size = 10

# --- paragraph 24:   Summary: No mycobacteria seen. ---
$cell.Range.Paragraphs(24).Range.InsertParagraphAfter()
$p24 = $cell.Range.Paragraphs(25).Range
$p24.InsertAfter("  Summary: No mycobacteria seen.")
$p24 = $cell.Range.Paragraphs(25).Range
$p24.Font.Name = "Times New Roman"
$p24.Font.Color = 16711680
$p24.Font.Size = 10

# --- paragraph 25: 07/12 – UNITS/ML – **Positive** ---
$cell.Range.Paragraphs(25).Range.InsertParagraphAfter()
$p25 = $cell.Range.Paragraphs(26).Range
$p25.InsertAfter("07/12 – UNITS/ML – **Positive**")
$p25 = $cell.Range.Paragraphs(26).Range
$p25.Font.Name = "Times New Roman"
$p25.Font.Color = 16711680
$p25.Font.Size = 10

# --- paragraph 26:   Summary: 2.94Log(10)IU/mL viral load detected. ---
$cell.Range.Paragraphs(26).Range.InsertParagraphAfter()
$p26 = $cell.Range.Paragraphs(27).Range
$p26.InsertAfter("  Summary: 2.94Log(10)IU/mL viral load detected.")
$p26 = $cell.Range.Paragraphs(27).Range
$p26.Font.Name = "Times New Roman"
$p26.Font.Color = 16711680
$p26.Font.Size = 10

# --- paragraph 27: 07/12 – SYPHILIS/TREPONEMA AB – Negative ---
$cell.Range.Paragraphs(27).Range.InsertParagraphAfter()
$p27 = $cell.Range.Paragraphs(28).Range
$p27.InsertAfter("07/12 – SYPHILIS/TREPONEMA AB – Negative")
$p27 = $cell.Range.Paragraphs(28).Range
$p27.Font.Name = "Times New Roman"
$p27.Font.Color = 16711680
$p27.Font.Size = 10

# --- paragraph 28: 07/12 – CMV DNA – Positive ---
$cell.Range.Paragraphs(28).Range.InsertParagraphAfter()
$p28 = $cell.Range.Paragraphs(29).Range
$p28.InsertAfter("07/12 – CMV DNA – Positive")
$p28 = $cell.Range.Paragraphs(29).Range
$p28.Font.Name = "Times New Roman"
$p28.Font.Color = 16711680
$p28.Font.Size = 10

# --- paragraph 29: 06/12 – BLC – PICC LINE NO GROWTH AFTER 5 DAYS ---
$cell.Range.Paragraphs(29).Range.InsertParagraphAfter()
$p29 = $cell.Range.Paragraphs(30).Range
$p29.InsertAfter("06/12 – BLC – PICC LINE NO GROWTH AFTER 5 DAYS")
$p29 = $cell.Range.Paragraphs(30).Range
$p29.Font.Name = "Times New Roman"
$p29.Font.Color = 16711680
$p29.Font.Size = 10

# --- paragraph 30: --------Previous result (1 year)-------- ---
$cell.Range.Paragraphs(30).Range.InsertParagraphAfter()
$p30 = $cell.Range.Paragraphs(31).Range
$p30.InsertAfter("--------Previous result (1 year)--------")
$p30 = $cell.Range.Paragraphs(31).Range
$p30.Font.Name = "Times New Roman"
$p30.Font.Color = 16711680
$p30.Font.Size = 10

# --- paragraph 31: 30/08 – CMV DNA – Positive ---
$cell.Range.Paragraphs(31).Range.InsertParagraphAfter()
$p31 = $cell.Range.Paragraphs(32).Range
$p31.InsertAfter("30/08 – CMV DNA – Positive")
$p31 = $cell.Range.Paragraphs(32).Range
$p31.Font.Name = "Times New Roman"
$p31.Font.Color = 16711680
$p31.Font.Size = 10

# --- paragraph 32: 19/08 – COMMENTS: – **No clear Result** + **Summar ---
$cell.Range.Paragraphs(32).Range.InsertParagraphAfter()
$p32 = $cell.Range.Paragraphs(33).Range
$p32.InsertAfter("19/08 – COMMENTS: – **No clear Result** + **Summary:** Not tested for toxoplasma (IgG negative).  ")
$p32 = $cell.Range.Paragraphs(33).Range
$p32.Font.Name = "Times New Roman"
$p32.Font.Color = 16711680
$p32.Font.Size = 10

# --- paragraph 33: (empty) ---
$cell.Range.Paragraphs(33).Range.InsertParagraphAfter()
$p33 = $cell.Range.Paragraphs(34).Range
$p33.Font.Name = "Times New Roman"
$p33.Font.Color = 16711680
$p33.Font.Size = 10

# --- paragraph 34: **Explanation**: The microbiology test for toxopla ---
$cell.Range.Paragraphs(34).Range.InsertParagraphAfter()
$p34 = $cell.Range.Paragraphs(35).Range
$p34.InsertAfter("**Explanation**: The microbiology test for toxoplasma was not performed due to a prior negative IgG result in serum. Since no direct microbiological analysis (e.g., culture, PCR) for toxoplasma was conducted, there is no definitive pathogen–specific finding to classify as Positive/Negative.")
$p34 = $cell.Range.Paragraphs(35).Range
$p34.Font.Name = "Times New Roman"
$p34.Font.Color = 16711680
$p34.Font.Size = 10

# --- paragraph 35: 08/08 – TB CULTURE – **Positive**   ---
$cell.Range.Paragraphs(35).Range.InsertParagraphAfter()
$p35 = $cell.Range.Paragraphs(36).Range
$p35.InsertAfter("08/08 – TB CULTURE – **Positive**  ")
$p35 = $cell.Range.Paragraphs(36).Range
$p35.Font.Name = "Times New Roman"
$p35.Font.Color = 16711680
$p35.Font.Size = 10

# --- paragraph 36: **Summary:** Acid–fast bacillus isolated (previous ---
$cell.Range.Paragraphs(36).Range.InsertParagraphAfter()
$p36 = $cell.Range.Paragraphs(37).Range
$p36.InsertAfter("**Summary:** Acid–fast bacillus isolated (previously noted).")
$p36 = $cell.Range.Paragraphs(37).Range
$p36.Font.Name = "Times New Roman"
$p36.Font.Color = 16711680
$p36.Font.Size = 10

# --- paragraph 37: 08/08 – CMV DNA – Positive ---
$cell.Range.Paragraphs(37).Range.InsertParagraphAfter()
$p37 = $cell.Range.Paragraphs(38).Range
$p37.InsertAfter("08/08 – CMV DNA – Positive")
$p37 = $cell.Range.Paragraphs(38).Range
$p37.Font.Name = "Times New Roman"
$p37.Font.Color = 16711680
$p37.Font.Size = 10

# --- paragraph 38: 08/08 – FUNGUS CULTURE – **Positive**   ---
$cell.Range.Paragraphs(38).Range.InsertParagraphAfter()
$p38 = $cell.Range.Paragraphs(39).Range
$p38.InsertAfter("08/08 – FUNGUS CULTURE – **Positive**  ")
$p38 = $cell.Range.Paragraphs(39).Range
$p38.Font.Name = "Times New Roman"
$p38.Font.Color = 16711680
$p38.Font.Size = 10

# --- paragraph 39: **Summary:** Candida albicans isolated. ---
$cell.Range.Paragraphs(39).Range.InsertParagraphAfter()
$p39 = $cell.Range.Paragraphs(40).Range
$p39.InsertAfter("**Summary:** Candida albicans isolated.")
$p39 = $cell.Range.Paragraphs(40).Range
$p39.Font.Name = "Times New Roman"
$p39.Font.Color = 16711680
$p39.Font.Size = 10

# --- paragraph 40: 08/08 – EBV VCA IgG – Positive ---
$cell.Range.Paragraphs(40).Range.InsertParagraphAfter()
$p40 = $cell.Range.Paragraphs(41).Range
$p40.InsertAfter("08/08 – EBV VCA IgG – Positive")
$p40 = $cell.Range.Paragraphs(41).Range
$p40.Font.Name = "Times New Roman"
$p40.Font.Color = 16711680
$p40.Font.Size = 10

# --- paragraph 41: 08/08 – RESP. CULT AND MICRO – **Positive**   ---
$cell.Range.Paragraphs(41).Range.InsertParagraphAfter()
$p41 = $cell.Range.Paragraphs(42).Range
$p41.InsertAfter("08/08 – RESP. CULT AND MICRO – **Positive**  ")
$p41 = $cell.Range.Paragraphs(42).Range
$p41.Font.Name = "Times New Roman"
$p41.Font.Color = 16711680
$p41.Font.Size = 10

# --- paragraph 42: **Summary:** Streptococcus viridans 10^4 – 10^5 or ---
$cell.Range.Paragraphs(42).Range.InsertParagraphAfter()
$p42 = $cell.Range.Paragraphs(43).Range
$p42.InsertAfter("**Summary:** Streptococcus viridans 10^4 – 10^5 orgs/mL.")
$p42 = $cell.Range.Paragraphs(43).Range
$p42.Font.Name = "Times New Roman"
$p42.Font.Color = 16711680
$p42.Font.Size = 10

# --- paragraph 43: 07/08 – URINE CULTURE – **Positive**   ---
$cell.Range.Paragraphs(43).Range.InsertParagraphAfter()
$p43 = $cell.Range.Paragraphs(44).Range
$p43.InsertAfter("07/08 – URINE CULTURE – **Positive**  ")
$p43 = $cell.Range.Paragraphs(44).Range
$p43.Font.Name = "Times New Roman"
$p43.Font.Color = 16711680
$p43.Font.Size = 10

# --- paragraph 44: **Summary:** Candida species isolated (10^4 – 10^5 ---
$cell.Range.Paragraphs(44).Range.InsertParagraphAfter()
$p44 = $cell.Range.Paragraphs(45).Range
$p44.InsertAfter("**Summary:** Candida species isolated (10^4 – 10^5 CFU/mL).")
$p44 = $cell.Range.Paragraphs(45).Range
$p44.Font.Name = "Times New Roman"
$p44.Font.Color = 16711680
$p44.Font.Size = 10

# --- paragraph 45: 06/08 – HEPATITIS A IgG – Positive ---
$cell.Range.Paragraphs(45).Range.InsertParagraphAfter()
$p45 = $cell.Range.Paragraphs(46).Range
$p45.InsertAfter("06/08 – HEPATITIS A IgG – Positive")
$p45 = $cell.Range.Paragraphs(46).Range
$p45.Font.Name = "Times New Roman"
$p45.Font.Color = 16711680
$p45.Font.Size = 10

# --- paragraph 46: 06/08 – VZV IgG – Positive ---
$cell.Range.Paragraphs(46).Range.InsertParagraphAfter()
$p46 = $cell.Range.Paragraphs(47).Range
$p46.InsertAfter("06/08 – VZV IgG – Positive")
$p46 = $cell.Range.Paragraphs(47).Range
$p46.Font.Name = "Times New Roman"
$p46.Font.Color = 16711680
$p46.Font.Size = 10

# --- paragraph 47: 06/08 – CMV IgG – Positive ---
$cell.Range.Paragraphs(47).Range.InsertParagraphAfter()
$p47 = $cell.Range.Paragraphs(48).Range
$p47.InsertAfter("06/08 – CMV IgG – Positive")
$p47 = $cell.Range.Paragraphs(48).Range
$p47.Font.Name = "Times New Roman"
$p47.Font.Color = 16711680
$p47.Font.Size = 10

# --- paragraph 48: 06/08 – PCR FOR MTB COMPLEX – **Positive**   ---
$cell.Range.Paragraphs(48).Range.InsertParagraphAfter()
$p48 = $cell.Range.Paragraphs(49).Range
$p48.InsertAfter("06/08 – PCR FOR MTB COMPLEX – **Positive**  ")
$p48 = $cell.Range.Paragraphs(49).Range
$p48.Font.Name = "Times New Roman"
$p48.Font.Color = 16711680
$p48.Font.Size = 10

# --- paragraph 49: **Summary:** MTB complex detected, no Rifampicin r ---
$cell.Range.Paragraphs(49).Range.InsertParagraphAfter()
$p49 = $cell.Range.Paragraphs(50).Range
$p49.InsertAfter("**Summary:** MTB complex detected, no Rifampicin resistance.")
$p49 = $cell.Range.Paragraphs(50).Range
$p49.Font.Name = "Times New Roman"
$p49.Font.Color = 16711680
$p49.Font.Size = 10

# --- paragraph 50: 06/08 – TB CULTURE – **Positive**   ---
$cell.Range.Paragraphs(50).Range.InsertParagraphAfter()
$p50 = $cell.Range.Paragraphs(51).Range
$p50.InsertAfter("06/08 – TB CULTURE – **Positive**  ")
$p50 = $cell.Range.Paragraphs(51).Range
$p50.Font.Name = "Times New Roman"
$p50.Font.Color = 16711680
$p50.Font.Size = 10

# --- paragraph 51: **Summary:** Mycobacterium tuberculosis complex id ---
$cell.Range.Paragraphs(51).Range.InsertParagraphAfter()
$p51 = $cell.Range.Paragraphs(52).Range
$p51.InsertAfter("**Summary:** Mycobacterium tuberculosis complex identified.")
$p51 = $cell.Range.Paragraphs(52).Range
$p51.Font.Name = "Times New Roman"
$p51.Font.Color = 16711680
$p51.Font.Size = 10

# --- paragraph 52: 06/08 – PCR FOR MTB COMPLEX – **Positive**   ---
$cell.Range.Paragraphs(52).Range.InsertParagraphAfter()
$p52 = $cell.Range.Paragraphs(53).Range
$p52.InsertAfter("06/08 – PCR FOR MTB COMPLEX – **Positive**  ")
$p52 = $cell.Range.Paragraphs(53).Range
$p52.Font.Name = "Times New Roman"
$p52.Font.Color = 16711680
$p52.Font.Size = 10

# --- paragraph 53: **Summary:** MTB detected, no Rifampicin resistanc ---
$cell.Range.Paragraphs(53).Range.InsertParagraphAfter()
$p53 = $cell.Range.Paragraphs(54).Range
$p53.InsertAfter("**Summary:** MTB detected, no Rifampicin resistance.")
$p53 = $cell.Range.Paragraphs(54).Range
$p53.Font.Name = "Times New Roman"
$p53.Font.Color = 16711680
$p53.Font.Size = 10

# --- paragraph 54: 06/08 – TB CULTURE – **Positive**   ---
$cell.Range.Paragraphs(54).Range.InsertParagraphAfter()
$p54 = $cell.Range.Paragraphs(55).Range
$p54.InsertAfter("06/08 – TB CULTURE – **Positive**  ")
$p54 = $cell.Range.Paragraphs(55).Range
$p54.Font.Name = "Times New Roman"
$p54.Font.Color = 16711680
$p54.Font.Size = 10

# --- paragraph 55: **Summary:** Acid–alcohol fast bacillus isolated. ---
$cell.Range.Paragraphs(55).Range.InsertParagraphAfter()
$p55 = $cell.Range.Paragraphs(56).Range
$p55.InsertAfter("**Summary:** Acid–alcohol fast bacillus isolated.")
$p55 = $cell.Range.Paragraphs(56).Range
$p55.Font.Name = "Times New Roman"
$p55.Font.Color = 16711680
$p55.Font.Size = 10

# --- paragraph 56: 06/08 – HIV 1 RNA – Positive ---
$cell.Range.Paragraphs(56).Range.InsertParagraphAfter()
$p56 = $cell.Range.Paragraphs(57).Range
$p56.InsertAfter("06/08 – HIV 1 RNA – Positive")
$p56 = $cell.Range.Paragraphs(57).Range
$p56.Font.Name = "Times New Roman"
$p56.Font.Color = 16711680
$p56.Font.Size = 10

# --- paragraph 57: 06/08 – TB CULTURE – **Positive**   ---
$cell.Range.Paragraphs(57).Range.InsertParagraphAfter()
$p57 = $cell.Range.Paragraphs(58).Range
$p57.InsertAfter("06/08 – TB CULTURE – **Positive**  ")
$p57 = $cell.Range.Paragraphs(58).Range
$p57.Font.Name = "Times New Roman"
$p57.Font.Color = 16711680
$p57.Font.Size = 10

# --- paragraph 58: **Summary:** Mycobacterium tuberculosis complex id ---
$cell.Range.Paragraphs(58).Range.InsertParagraphAfter()
$p58 = $cell.Range.Paragraphs(59).Range
$p58.InsertAfter("**Summary:** Mycobacterium tuberculosis complex identified.")
$p58 = $cell.Range.Paragraphs(59).Range
$p58.Font.Name = "Times New Roman"
$p58.Font.Color = 16711680
$p58.Font.Size = 10

# --- paragraph 59: 06/08 – COPIES/ML – **Positive**   ---
$cell.Range.Paragraphs(59).Range.InsertParagraphAfter()
$p59 = $cell.Range.Paragraphs(60).Range
$p59.InsertAfter("06/08 – COPIES/ML – **Positive**  ")
$p59 = $cell.Range.Paragraphs(60).Range
$p59.Font.Name = "Times New Roman"
$p59.Font.Color = 16711680
$p59.Font.Size = 10

# --- paragraph 60: **Summary:** 6.47Log(10)copies/mL HIV viral load. ---
$cell.Range.Paragraphs(60).Range.InsertParagraphAfter()
$p60 = $cell.Range.Paragraphs(61).Range
$p60.InsertAfter("**Summary:** 6.47Log(10)copies/mL HIV viral load.")
$p60 = $cell.Range.Paragraphs(61).Range
$p60.Font.Name = "Times New Roman"
$p60.Font.Color = 16711680
$p60.Font.Size = 10

Write-Host "Micro results section rebuilt"